# Remove duplicate rows that were accidentally included twice.
# Row 122 duplicated row 121 (2_06Nov23 / F_sed_4m_R2_F1_06Nov23 / F_sed_4m_R2_F2_06Nov23)
# Row 131 duplicated row 130 (2_17Nov23 / F_sed_4m_R2_F1_17Nov23 / F_sed_4m_R2_F2_17Nov23)
# Deleting the higher-numbered row first keeps the lower row index valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(131).Delete()
$ws.Rows.Item(122).Delete()

# The row delete above can leave the running row-counter helper formulas in
# columns D and F ("=D(row-1)+1" / "=F(row-1)+1") broken for every row that
# shifted. Re-stamp them so they keep counting cleanly to the new last row.
$lastRow = $ws.UsedRange.Rows.Count
For ($r = 3; $r -le $lastRow; $r++) {
    $prev = $r - 1
    $ws.Range("D" + $r).Formula = "=D" + $prev + "+1"
    $ws.Range("F" + $r).Formula = "=F" + $prev + "+1"
}

$ws.Range("H135").Select()
